$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 11112247
$ws.Range("I15").Value = 11112247
$ws.Range("K15").Value = 33336741
$ws.Range("M15").Value = -33336572
$ws.Range("H51").Value = 5207.846
$ws.Range("J51").Value = 5036
$ws.Range("L51").Value = 5036
$ws.Range("N51").Value = -6004
$ws.Range("H92").Value = 1107.7142
$ws.Range("I92").Value = 959
$ws.Range("K92").Value = 959
$ws.Range("M92").Value = 289
$ws.Range("H100").Value = 1398.3334
$ws.Range("J100").Value = 806
$ws.Range("L100").Value = 806
$ws.Range("N100").Value = -1888
$ws.Range("H135").Value = 5833.222
$ws.Range("I135").Value = 1591
$ws.Range("K135").Value = 14319
$ws.Range("M135").Value = -11784
$ws.Range("H138").Value = 5163.8945
$ws.Range("I138").Value = 1978.5385
$ws.Range("K138").Value = 5935.6155
$ws.Range("M138").Value = -795.6154999999999
$ws.Range("H141").Value = 6699.778
$ws.Range("J141").Value = 4449.5
$ws.Range("L141").Value = 13348.5
$ws.Range("N141").Value = -23708.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16939.465
$ws.Range("I32").Value = 18516.65
$ws.Range("K32").Value = 18516.65
$ws.Range("M32").Value = -18229.65
$ws.Range("H45").Value = 3711.6206
$ws.Range("I45").Value = 3569.95
$ws.Range("J45").Value = 4026.4443
$ws.Range("K45").Value = 3569.95
$ws.Range("L45").Value = 4026.4443
$ws.Range("M45").Value = -3192.95
$ws.Range("N45").Value = -4780.4443
$ws.Range("H61").Value = 12104.3
$ws.Range("I61").Value = 13552.235
$ws.Range("K61").Value = 13552.235
$ws.Range("M61").Value = -13340.235
$ws.Range("H88").Value = 126874.75
$ws.Range("J88").Value = 126874.75
$ws.Range("L88").Value = 126874.75
$ws.Range("N88").Value = -127686.75
$ws.Range("H91").Value = 126874.75
$ws.Range("J91").Value = 126874.75
$ws.Range("L91").Value = 126874.75
$ws.Range("N91").Value = -129682.75
$ws.Range("H110").Value = 1706468.8
$ws.Range("I110").Value = 4085234.5
$ws.Range("J110").Value = 7350.4287
$ws.Range("K110").Value = 4085234.5
$ws.Range("L110").Value = 7350.4287
$ws.Range("M110").Value = -4083189.5
$ws.Range("N110").Value = -11440.4287
$ws.Range("H132").Value = 12697.188
$ws.Range("I132").Value = 14622.667
$ws.Range("J132").Value = 5345.364
$ws.Range("K132").Value = 43868.001
$ws.Range("L132").Value = 16036.092
$ws.Range("M132").Value = -41338.001
$ws.Range("N132").Value = -21096.092
$ws.Range("H136").Value = 12104.3
$ws.Range("I136").Value = 13552.235
$ws.Range("K136").Value = 40656.705
$ws.Range("M136").Value = -38106.705

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1411.0714
$ws.Range("I86").Value = 1535.6
$ws.Range("J86").Value = 1099.75
$ws.Range("K86").Value = 1535.6
$ws.Range("L86").Value = 1099.75
$ws.Range("M86").Value = -412.5999999999999
$ws.Range("N86").Value = -3345.75
$ws.Range("H89").Value = 1411.0714
$ws.Range("I89").Value = 1535.6
$ws.Range("J89").Value = 1099.75
$ws.Range("K89").Value = 7678
$ws.Range("L89").Value = 5498.75
$ws.Range("M89").Value = -2062
$ws.Range("N89").Value = -16730.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 225.8
$ws.Range("J7").Value = 255.8
$ws.Range("L7").Value = 255.8
$ws.Range("N7").Value = -481.8
$ws.Range("H132").Value = 47650588
$ws.Range("I132").Value = 83386280
$ws.Range("J132").Value = 2999.3333
$ws.Range("K132").Value = 250158840
$ws.Range("L132").Value = 8997.999899999999
$ws.Range("M132").Value = -250156310
$ws.Range("N132").Value = -14057.9999
$ws.Range("H133").Value = 68529
$ws.Range("J133").Value = 68529
$ws.Range("L133").Value = 68529
$ws.Range("N133").Value = -73589

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 179.88889
$ws.Range("I12").Value = 30.666666
$ws.Range("J12").Value = 254.5
$ws.Range("K12").Value = 91.99999800000001
$ws.Range("L12").Value = 763.5
$ws.Range("M12").Value = 81.00000199999999
$ws.Range("N12").Value = -1109.5
$ws.Range("H23").Value = 100000300
$ws.Range("J23").Value = 125000350
$ws.Range("L23").Value = 375001050
$ws.Range("N23").Value = -375001520
$ws.Range("H33").Value = 233
$ws.Range("J33").Value = 163.375
$ws.Range("L33").Value = 980.25
$ws.Range("N33").Value = -1546.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 437.8889
$ws.Range("I22").Value = 298.5
$ws.Range("J22").Value = 716.6667
$ws.Range("K22").Value = 298.5
$ws.Range("L22").Value = 716.6667
$ws.Range("M22").Value = -3.5
$ws.Range("N22").Value = -1306.6667
$ws.Range("H27").Value = 437.8889
$ws.Range("I27").Value = 298.5
$ws.Range("J27").Value = 716.6667
$ws.Range("K27").Value = 298.5
$ws.Range("L27").Value = 716.6667
$ws.Range("M27").Value = -191.5
$ws.Range("N27").Value = -930.6667
$ws.Range("H68").Value = 1895429.8
$ws.Range("I68").Value = 2067438.5
$ws.Range("K68").Value = 2067438.5
$ws.Range("M68").Value = -2066689.5
$ws.Range("H71").Value = 1895429.8
$ws.Range("I71").Value = 2067438.5
$ws.Range("K71").Value = 10337192.5
$ws.Range("M71").Value = -10333448.5
$ws.Range("H82").Value = 7813600
$ws.Range("I82").Value = 15625474
$ws.Range("K82").Value = 15625474
$ws.Range("M82").Value = -15625113
$ws.Range("H85").Value = 7813600
$ws.Range("I85").Value = 15625474
$ws.Range("K85").Value = 15625474
$ws.Range("M85").Value = -15624226
$ws.Range("H132").Value = 3830.7754
$ws.Range("I132").Value = 2995.5278
$ws.Range("K132").Value = 8986.5834
$ws.Range("M132").Value = -6456.5834
$ws.Range("H136").Value = 2745.3333
$ws.Range("I136").Value = 2023.7727
$ws.Range("J136").Value = 5920.2
$ws.Range("K136").Value = 6071.3181
$ws.Range("L136").Value = 17760.6
$ws.Range("M136").Value = -3521.3181
$ws.Range("N136").Value = -22860.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3027.36
$ws.Range("I107").Value = 3128.2942
$ws.Range("J107").Value = 2812.875
$ws.Range("K107").Value = 9384.882599999999
$ws.Range("L107").Value = 8438.625
$ws.Range("M107").Value = -7464.882599999999
$ws.Range("N107").Value = -12278.625
$ws.Range("H132").Value = 24160060
$ws.Range("J132").Value = 38470410
$ws.Range("L132").Value = 115411230
$ws.Range("N132").Value = -115416290
